# Update column F ("dSF") values on the data sheet with newly repulled data.
# The sheet layout: A=idx, B=date, C=TB, D=PC, E=dS0, F=dSF, G=K, H=IP, I=I0, J=IF
# Only column F changes; each row gets a fresh numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = 3
    4  = 7
    5  = 3
    8  = 3
    10 = 2
    11 = 6
    13 = 1
    14 = 3
    15 = 1
    16 = 4
    19 = 1
    20 = 2
    23 = -3
    24 = 2
    25 = -1
    28 = -4
    29 = -3
    33 = -4
    34 = 4
    36 = -3
    37 = -5
    38 = -8
    44 = -6
    45 = -2
    46 = -3
    47 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
